# Avances Etiquetado Roboflow 6/18/2025
# Add the new day's row of progress data (row 39) and update the
# active view/selection to reflect where the user was working.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New entry for 18/6/2025
$ws.Range("D39").Value = "18/6/2025"
$ws.Range("E39").Value = 373
$ws.Range("F39").Value = 553
$ws.Range("G39").Value = 0
$ws.Range("H39").Value = 0
$ws.Range("I39").Value = 1012
$ws.Range("J39").Value = "N/A"

# Update the visible scroll position / selection as recorded when the
# workbook was saved (topLeftCell A22, selection H41).
$win = $excel.ActiveWindow
$win.ScrollRow = 22
$win.ScrollColumn = 1
$ws.Range("H41").Select()
